$p = $ppt.ActivePresentation

# --- 1. Update the cached "date updates automatically" field text
#     (10 slide layouts + the slide master) from 7/23/20 -> 6/17/21
function Set-DatePlaceholderText($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "6/17/21"
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes
}

# --- 2. Fix up the braille-label text boxes on slide 2 (the braille rock-cycle slide)
$s = $p.Slides.Item(2)

$s.Shapes.Item("TextBox 3").TextFrame.TextRange.Text = "⠺⠂⠮⠗⠬⠀⠯⠻⠕⠨⠝"
$s.Shapes.Item("TextBox 16").TextFrame.TextRange.Text = "⠎⠫⠊;t⠎"
$s.Shapes.Item("TextBox 17").TextFrame.TextRange.Text = "⠎⠫⠊;t⠜⠽⠀⠗⠕⠉⠅⠎"
$s.Shapes.Item("TextBox 21").TextFrame.TextRange.Text = "⠍⠑⠞⠁⠍⠕⠗⠏⠓⠊⠉⠗⠕⠉⠅⠎"
$s.Shapes.Item("TextBox 23").TextFrame.TextRange.Text = "⠊⠛⠝⠑⠳⠎⠗⠕⠉⠅⠎"
$s.Shapes.Item("Rectangle 6").TextFrame.TextRange.Text = "⠠! ⠠⠗⠕⠉⠅⠀⠠⠉⠽⠉⠇⠑"
$s.Shapes.Item("Rectangle 38").TextFrame.TextRange.Text = "⠓⠂⠞⠀⠯⠀ ⠏⠗⠑⠎⠎⠥⠗⠑"

# Widen the title textbox so the longer braille title still fits
$s.Shapes.Item("Rectangle 6").Width = 3526928 / 12700.0

Write-Output "edit complete"
